$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1666666666666667
$ws.Range("C2").Value = 0.5894308943089431
$ws.Range("J2").Value = 0.01219512195121951
$ws.Range("P2").Value = 0.1219512195121951
$ws.Range("S2").Value = 0.1097560975609756
$ws.Range("B3").Value = 0.01744186046511628
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.04651162790697674
$ws.Range("P3").Value = 0.7034883720930233
$ws.Range("S3").Value = 0.2093023255813954
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("O4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.6222222222222222
$ws.Range("S4").Value = 0.3111111111111111
$ws.Range("B6").Value = 0.09042553191489362
$ws.Range("D6").Value = 0.005319148936170213
$ws.Range("F6").Value = 0.09574468085106383
$ws.Range("J6").Value = 0.2553191489361702
$ws.Range("O6").Value = 0.03191489361702127
$ws.Range("Q6").Value = 0.1276595744680851
$ws.Range("R6").Value = 0.04787234042553191
$ws.Range("S6").Value = 0.3457446808510639
$ws.Range("B7").Value = 0.1055555555555556
$ws.Range("D7").Value = 0.02777777777777778
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1166666666666667
$ws.Range("O7").Value = 0.02222222222222222
$ws.Range("Q7").Value = 0.2055555555555555
$ws.Range("R7").Value = 0.08888888888888889
$ws.Range("S7").Value = 0.3666666666666666
$ws.Range("B8").Value = 0.07942973523421588
$ws.Range("D8").Value = 0.0285132382892057
$ws.Range("F8").Value = 0.06924643584521385
$ws.Range("J8").Value = 0.1181262729124236
$ws.Range("O8").Value = 0.01629327902240326
$ws.Range("Q8").Value = 0.1710794297352342
$ws.Range("R8").Value = 0.06109979633401222
$ws.Range("S8").Value = 0.4562118126272913
$ws.Range("B9").Value = 0.1013513513513514
$ws.Range("D9").Value = 0.01351351351351351
$ws.Range("F9").Value = 0.03378378378378379
$ws.Range("J9").Value = 0.1554054054054054
$ws.Range("O9").Value = 0.01351351351351351
$ws.Range("Q9").Value = 0.1486486486486487
$ws.Range("R9").Value = 0.08783783783783784
$ws.Range("S9").Value = 0.4459459459459459
$ws.Range("B10").Value = 0.09811320754716982
$ws.Range("D10").Value = 0.02264150943396226
$ws.Range("E10").Value = 0.0009433962264150943
$ws.Range("F10").Value = 0.05849056603773585
$ws.Range("J10").Value = 0.1245283018867925
$ws.Range("O10").Value = 0.02075471698113207
$ws.Range("Q10").Value = 0.230188679245283
$ws.Range("R10").Value = 0.07264150943396226
$ws.Range("S10").Value = 0.3716981132075471
$ws.Range("F11").Value = 0.003472222222222222
$ws.Range("G11").Value = 0.1423611111111111
$ws.Range("J11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.2048611111111111
$ws.Range("L11").Value = 0.5104166666666666
$ws.Range("S11").Value = 0.02777777777777778
$ws.Range("G12").Value = 0.7225806451612903
$ws.Range("J12").Value = 0.1870967741935484
$ws.Range("K12").Value = 0.02580645161290323
$ws.Range("L12").Value = 0.03870967741935484
$ws.Range("S12").Value = 0.02580645161290323
$ws.Range("G13").Value = 0.813953488372093
$ws.Range("J13").Value = 0.1627906976744186
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.02538071065989848
$ws.Range("H15").Value = 0.1624365482233502
$ws.Range("I15").Value = 0.07106598984771574
$ws.Range("J15").Value = 0.3401015228426396
$ws.Range("K15").Value = 0.07614213197969544
$ws.Range("M15").Value = 0.005076142131979695
$ws.Range("O15").Value = 0.06091370558375635
$ws.Range("S15").Value = 0.2588832487309645
$ws.Range("F16").Value = 0.02958579881656805
$ws.Range("H16").Value = 0.1952662721893491
$ws.Range("I16").Value = 0.04733727810650887
$ws.Range("J16").Value = 0.3727810650887574
$ws.Range("K16").Value = 0.1183431952662722
$ws.Range("M16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.07100591715976332
$ws.Range("S16").Value = 0.1597633136094675
$ws.Range("F17").Value = 0.02870813397129187
$ws.Range("H17").Value = 0.1770334928229665
$ws.Range("I17").Value = 0.1124401913875598
$ws.Range("J17").Value = 0.3875598086124402
$ws.Range("K17").Value = 0.08851674641148326
$ws.Range("M17").Value = 0.0215311004784689
$ws.Range("O17").Value = 0.05263157894736842
$ws.Range("S17").Value = 0.131578947368421
$ws.Range("F18").Value = 0.01342281879194631
$ws.Range("H18").Value = 0.1812080536912752
$ws.Range("I18").Value = 0.08053691275167785
$ws.Range("J18").Value = 0.3221476510067114
$ws.Range("K18").Value = 0.1208053691275168
$ws.Range("M18").Value = 0.02013422818791946
$ws.Range("O18").Value = 0.04697986577181208
$ws.Range("S18").Value = 0.2147651006711409
$ws.Range("F19").Value = 0.01655629139072848
$ws.Range("H19").Value = 0.2682119205298013
$ws.Range("I19").Value = 0.05463576158940397
$ws.Range("J19").Value = 0.3294701986754967
$ws.Range("K19").Value = 0.1142384105960265
$ws.Range("M19").Value = 0.02400662251655629
$ws.Range("N19").Value = 0.001655629139072848
$ws.Range("O19").Value = 0.05960264900662252
$ws.Range("S19").Value = 0.1316225165562914
